$wb = $excel.ActiveWorkbook

# --- Update the "Date" metadata value on the Metadata sheet ---
$metaSheet = $wb.Worksheets.Item("Metadata")
$metaSheet.Range("B8").Value = "2025-09-23T20:31:36+00:00"

# --- Fix the "Adult" concept code on the Concepts sheet ---
# B5 previously held "Pediatric-and-Adult" (code) while C5 holds
# "Pediatric and Adult" (display). The code should match the display text.
$conceptsSheet = $wb.Worksheets.Item("Concepts")
$conceptsSheet.Range("B5").Value = "Pediatric and Adult"
